$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.Value = '''62.953.84'
$cell.Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -1.19%  '

$cell = $ws.Cells.Item(3, 4)
$cell.Value = '''3.162.69'
$cell.Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  +1.36%  '

$ws.Cells.Item(4, 5).Value = '  -0.05%  '

$cell = $ws.Cells.Item(5, 4)
$cell.Value = '''587.67'
$cell.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -1.74%  '

$cell = $ws.Cells.Item(6, 4)
$cell.Value = '''138.45'
$cell.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -2.47%  '

$ws.Cells.Item(7, 5).Value = '  -0.03%  '

$cell = $ws.Cells.Item(8, 4)
$cell.Value = '''3.158.02'
$cell.Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +1.34%  '

$ws.Cells.Item(9, 5).Value = '  -1.15%  '

$ws.Cells.Item(10, 5).Value = '  -1.30%  '

$cell = $ws.Cells.Item(11, 4)
$cell.Value = '''5.30'
$cell.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -0.63%  '

$cell = $ws.Cells.Item(12, 4)
$cell.Value = '''0.458'
$cell.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -1.29%  '

$cell = $ws.Cells.Item(13, 4)
$cell.Value = '''0.0000245'
$cell.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -3.13%  '

$cell = $ws.Cells.Item(14, 4)
$cell.Value = '''34.11'
$cell.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -2.38%  '

$cell = $ws.Cells.Item(15, 4)
$cell.Value = '''3.678.33'
$cell.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +1.12%  '

$ws.Cells.Item(16, 5).Value = '  +0.74%  '

$cell = $ws.Cells.Item(17, 4)
$cell.Value = '''3.158.08'
$cell.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +1.27%  '

$cell = $ws.Cells.Item(18, 4)
$cell.Value = '''62.930.42'
$cell.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -1.29%  '

$cell = $ws.Cells.Item(19, 4)
$cell.Value = '''6.66'
$cell.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -2.06%  '

$cell = $ws.Cells.Item(20, 4)
$cell.Value = '''477.31'
$cell.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -1.19%  '

$cell = $ws.Cells.Item(22, 4)
$cell.Value = '''0.701'
$cell.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.40%  '

$cell = $ws.Cells.Item(23, 4)
$cell.Value = '''7.74'
$cell.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +1.98%  '

$cell = $ws.Cells.Item(24, 4)
$cell.Value = '''84.37'
$cell.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -2.56%  '

$cell = $ws.Cells.Item(25, 4)
$cell.Value = '''12.98'
$cell.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -2.85%  '

$ws.Cells.Item(26, 5).Value = '  -0.01%  '

$ws.Cells.Item(27, 5).Value = '  -1.04%  '

$cell = $ws.Cells.Item(28, 4)
$cell.Value = '''7.11'
$cell.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +2.26%  '

$cell = $ws.Cells.Item(29, 4)
$cell.Value = '''7.94'
$cell.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -3.27%  '

$cell = $ws.Cells.Item(30, 4)
$cell.Value = '''2.08'
$cell.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +1.36%  '

$ws.Cells.Item(31, 5).Value = '  -0.05%  '

$cell = $ws.Cells.Item(32, 4)
$cell.Value = '''26.79'
$cell.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -0.86%  '

$cell = $ws.Cells.Item(33, 4)
$cell.Value = '''0.106'
$cell.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -4.62%  '

$ws.Cells.Item(34, 5).Value = '  -4.91%  '

$ws.Cells.Item(35, 5).Value = '  -3.13%  '

$cell = $ws.Cells.Item(36, 4)
$cell.Value = '''52.51'
$cell.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +0.09%  '

$cell = $ws.Cells.Item(37, 4)
$cell.Value = '''5.78'
$cell.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -2.91%  '

$cell = $ws.Cells.Item(38, 4)
$cell.Value = '''0.0₃0704'
$cell.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -4.79%  '

$ws.Cells.Item(39, 5).Value = '  -1.75%  '

$cell = $ws.Cells.Item(40, 4)
$cell.Value = '''416.95'
$cell.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -4.41%  '

$ws.Cells.Item(41, 5).Value = '  -5.99%  '

$cell = $ws.Cells.Item(42, 4)
$cell.Value = '''2.954.54'
$cell.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +3.06%  '

$cell = $ws.Cells.Item(43, 4)
$cell.Value = '''8.27'
$cell.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +0.46%  '

$ws.Cells.Item(44, 5).Value = '  -7.19%  '

$cell = $ws.Cells.Item(45, 4)
$cell.Value = '''0.259'
$cell.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +0.43%  '

$ws.Cells.Item(47, 5).Value = '  -3.15%  '

$cell = $ws.Cells.Item(48, 4)
$cell.Value = '''25.48'
$cell.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -1.15%  '

$ws.Cells.Item(49, 5).Value = '  -0.20%  '

$ws.Cells.Item(50, 5).Value = '  -4.26%  '

$cell = $ws.Cells.Item(51, 4)
$cell.Value = '''119.63'
$cell.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -1.62%  '
